$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 13 of portfolio data for 2025-08-28.
# The leading apostrophe forces Excel to treat the date-like string as
# literal text instead of auto-converting it into a date serial number,
# matching the existing rows (A2:A12) which are also stored as text.
$ws.Range("A13").Value = "'2025-08-28"
$ws.Range("A13").Style = "Normal"

$ws.Range("B13").Value = 56.31999969482422
$ws.Range("C13").Value = 675.4500122070312
$ws.Range("D13").Value = 315.5
